$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format ("@") on D-column cells whose new value would otherwise
# be auto-coerced to a number by Excel (single-dot decimal-looking strings),
# so they stay stored as text exactly like the source data.
$textCells = @("D4", "D5", "D6", "D10", "D11", "D12", "D19", "D20", "D23", "D24", "D26", "D28", "D29", "D32", "D34", "D35", "D37", "D41", "D42", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value = '50.132.08'
$ws.Range('E2').Value = '  +4.30%  '
$ws.Range('D3').Value = '2.661.49'
$ws.Range('E3').Value = '  +6.80%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '114.70'
$ws.Range('E5').Value = '  +8.70%  '
$ws.Range('D6').Value = '326.62'
$ws.Range('E6').Value = '  +2.83%  '
$ws.Range('E7').Value = '  +2.24%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  +3.75%  '
$ws.Range('D10').Value = '41.40'
$ws.Range('E10').Value = '  +6.45%  '
$ws.Range('D11').Value = '20.18'
$ws.Range('E11').Value = '  -0.18%  '
$ws.Range('D12').Value = '0.0827'
$ws.Range('E12').Value = '  +3.19%  '
$ws.Range('E13').Value = '  +0.45%  '
$ws.Range('E14').Value = '  +4.44%  '
$ws.Range('D15').Value = '3.075.11'
$ws.Range('E15').Value = '  +6.51%  '
$ws.Range('D16').Value = '2.644.79'
$ws.Range('E16').Value = '  +5.91%  '
$ws.Range('E17').Value = '  +5.98%  '
$ws.Range('D18').Value = '50.062.38'
$ws.Range('E18').Value = '  +4.31%  '
$ws.Range('D19').Value = '13.29'
$ws.Range('E19').Value = '  +4.00%  '
$ws.Range('D20').Value = '6.78'
$ws.Range('E20').Value = '  +3.07%  '
$ws.Range('E21').Value = '  -2.04%  '
$ws.Range('D22').Value = '0.0₃0961'
$ws.Range('E22').Value = '  +3.43%  '
$ws.Range('D23').Value = '72.55'
$ws.Range('E23').Value = '  +2.15%  '
$ws.Range('D24').Value = '276.70'
$ws.Range('E24').Value = '  +3.22%  '
$ws.Range('E25').Value = '  +3.92%  '
$ws.Range('D26').Value = '27.00'
$ws.Range('E26').Value = '  +5.00%  '
$ws.Range('E27').Value = '  +0.11%  '
$ws.Range('D28').Value = '10.08'
$ws.Range('E28').Value = '  +3.56%  '
$ws.Range('D29').Value = '37.01'
$ws.Range('E29').Value = '  +7.25%  '
$ws.Range('E30').Value = '  -2.14%  '
$ws.Range('E31').Value = '  +2.44%  '
$ws.Range('D32').Value = '50.23'
$ws.Range('E32').Value = '  +2.07%  '
$ws.Range('E33').Value = '  +4.32%  '
$ws.Range('D34').Value = '19.80'
$ws.Range('E34').Value = '  +3.89%  '
$ws.Range('D35').Value = '0.0816'
$ws.Range('D37').Value = '5.00'
$ws.Range('E37').Value = '  +9.31%  '
$ws.Range('E38').Value = '  +7.17%  '
$ws.Range('E39').Value = '  +9.05%  '
$ws.Range('E40').Value = '  +2.83%  '
$ws.Range('D41').Value = '124.44'
$ws.Range('E41').Value = '  +0.88%  '
$ws.Range('D42').Value = '22.50'
$ws.Range('E42').Value = '  +1.25%  '
$ws.Range('E43').Value = '  +0.09%  '
$ws.Range('E44').Value = '  +5.30%  '
$ws.Range('D45').Value = '2.107.32'
$ws.Range('E45').Value = '  +5.39%  '
$ws.Range('E46').Value = '  +5.92%  '
$ws.Range('E47').Value = '  +13.36%  '
$ws.Range('E48').Value = '  +4.69%  '
$ws.Range('E49').Value = '  +2.38%  '
$ws.Range('E50').Value = '  +3.65%  '
$ws.Range('D51').Value = '60.20'
$ws.Range('E51').Value = '  +6.56%  '
